{"js": "// \"updated to get the table references to not be bold.\"\n//\n// The document defines a custom paragraph style \"Table Caption\"\n// (styleId \"TableCaption\", based on the built-in \"Caption\" style) that is\n// applied to table-reference/caption paragraphs. Turn bold off on that\n// style so anything using it (directly or by inheritance) is no longer\n// rendered bold.\nconst tableCaptionStyle = context.document.getStyles().getByNameOrNullObject(\"Table Caption\");\ntableCaptionStyle.load(\"nameLocal\");\nawait context.sync();\n\nif (!tableCaptionStyle.isNullObject) {\n  tableCaptionStyle.font.bold = false;\n  await context.sync();\n}\n", "ps1": "# \"updated to get the table references to not be bold.\"\n#\n# The document defines a custom paragraph style \"Table Caption\" (based on\n# the built-in \"Caption\" style) that is applied to table-reference/caption\n# paragraphs. Turn bold off on that style so anything using it (directly or\n# by inheritance) is no longer rendered bold.\n$d = $word.ActiveDocument\n\n$tableCaptionStyle = $d.Styles(\"Table Caption\")\n$tableCaptionStyle.Font.Bold = 0\n"}
